$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New column F header and value
$ws.Range("F1").Value = "OpcionMotivo"
$ws.Range("F2").Value = "Anulación por desistimiento"

# Update NroPoliza value in E2 (kept as text, matching the existing quote-prefixed style)
$ws.Range("E2").Value = "'12112001837"

# Update selection to match target state
$ws.Range("F4").Select()
